$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''62.698.35'
$ws.Range("E2").Value = '''  +4.33%  '

# Row 3
$ws.Range("D3").Value = '''3.339.34'
$ws.Range("E3").Value = '''  +4.41%  '

# Row 4
$ws.Range("E4").Value = '''  -0.05%  '

# Row 5
$ws.Range("D5").Value = '''553.48'
$ws.Range("E5").Value = '''  +2.88%  '

# Row 6
$ws.Range("D6").Value = '''151.96'
$ws.Range("E6").Value = '''  +4.48%  '

# Row 7
$ws.Range("E7").Value = '''  -0.06%  '

# Row 8
$ws.Range("D8").Value = '''0.530'
$ws.Range("E8").Value = '''  +2.27%  '

# Row 9
$ws.Range("E9").Value = '''  +2.64%  '

# Row 10
$ws.Range("D10").Value = '''0.118'
$ws.Range("E10").Value = '''  +3.38%  '

# Row 11
$ws.Range("E11").Value = '''  +0.80%  '

# Row 12
$ws.Range("D12").Value = '''3.909.37'
$ws.Range("E12").Value = '''  +4.23%  '

# Row 13
$ws.Range("E13").Value = '''  -0.39%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''26.82'
$ws.Range("E14").Value = '''  +2.75%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000179'
$ws.Range("E15").Value = '''  +2.39%  '

# Row 16
$ws.Range("D16").Value = '''62.721.43'
$ws.Range("E16").Value = '''  +4.37%  '

# Row 17
$ws.Range("D17").Value = '''3.339.31'
$ws.Range("E17").Value = '''  +4.11%  '

# Row 18
$ws.Range("D18").Value = '''6.53'
$ws.Range("E18").Value = '''  +5.09%  '

# Row 19
$ws.Range("E19").Value = '''  +4.44%  '

# Row 20
$ws.Range("D20").Value = '''8.44'
$ws.Range("E20").Value = '''  +1.71%  '

# Row 21
$ws.Range("D21").Value = '''386.78'
$ws.Range("E21").Value = '''  -0.02%  '

# Row 22
$ws.Range("E22").Value = '''  +0.36%  '

# Row 23
$ws.Range("E23").Value = '''  +1.36%  '

# Row 24
$ws.Range("D24").Value = '''70.57'
$ws.Range("E24").Value = '''  +0.50%  '

# Row 25
$ws.Range("E25").Value = '''  +1.95%  '

# Row 26
$ws.Range("E26").Value = '''  -0.91%  '

# Row 27
$ws.Range("D27").Value = '''0.0₃0958'
$ws.Range("E27").Value = '''  +5.22%  '

# Row 28
$ws.Range("E28").Value = '''  +0.17%  '

# Row 29
$ws.Range("D29").Value = '''1.97'
$ws.Range("E29").Value = '''  +2.95%  '

# Row 30
$ws.Range("D30").Value = '''6.42'
$ws.Range("E30").Value = '''  +4.16%  '

# Row 31
$ws.Range("D31").Value = '''22.97'
$ws.Range("E31").Value = '''  +2.30%  '

# Row 32
$ws.Range("E32").Value = '''  +1.92%  '

# Row 33
$ws.Range("E33").Value = '''  +6.20%  '

# Row 34
$ws.Range("D34").Value = '''6.67'
$ws.Range("E34").Value = '''  +2.88%  '

# Row 35
$ws.Range("D35").Value = '''160.94'
$ws.Range("E35").Value = '''  +2.90%  '

# Row 36
$ws.Range("E36").Value = '''  +9.46%  '

# Row 37
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").Value = '''1.87'
$ws.Range("E37").Value = '''  +10.30%  '

# Row 38
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '''27.28'
$ws.Range("E38").Value = '''  +5.82%  '

# Row 39
$ws.Range("D39").Value = '''2.855.22'
$ws.Range("E39").Value = '''  +3.28%  '

# Row 40
$ws.Range("D40").Value = '''0.0737'
$ws.Range("E40").Value = '''  +3.21%  '

# Row 41
$ws.Range("D41").Value = '''0.0312'
$ws.Range("E41").Value = '''  +7.93%  '

# Row 42
$ws.Range("E42").Value = '''  +0.66%  '

# Row 43
$ws.Range("D43").Value = '''0.752'
$ws.Range("E43").Value = '''  +3.22%  '

# Row 44
$ws.Range("D44").Value = '''40.69'
$ws.Range("E44").Value = '''  +2.82%  '

# Row 45
$ws.Range("E45").Value = '''  +2.80%  '

# Row 46
$ws.Range("D46").Value = '''3.380.92'
$ws.Range("E46").Value = '''  +4.25%  '

# Row 47
$ws.Range("D47").Value = '''21.92'
$ws.Range("E47").Value = '''  +5.87%  '

# Row 48
$ws.Range("E48").Value = '''  +3.03%  '

# Row 49
$ws.Range("D49").Value = '''6.27'
$ws.Range("E49").Value = '''  +0.98%  '

# Row 50
$ws.Range("D50").Value = '''0.803'
$ws.Range("E50").Value = '''  +2.02%  '

# Row 51
$ws.Range("D51").Value = '''283.23'
$ws.Range("E51").Value = '''  +7.97%  '
